$d = $word.ActiveDocument

# --- Step 1: remove the hyperlink around "https://automationstepbystep.com/"
# but keep the plain text (no hyperlink, no character style) ---
$h = $d.Hyperlinks.Item(3)
$h.Delete()

$p3 = $d.Paragraphs.Item(3).Range
$pStart = $p3.Start
$pEnd = $p3.End

# Remove the old (styled) run's text, leaving the trailing field remnant
# character intact so the paragraph mark / pPr (ListParagraph + numPr) is
# preserved.
$oldText = $d.Range($pStart, $pEnd - 1)
$oldText.Delete()

# Insert a fresh, unformatted run with the same text right before that
# remnant character.
$insPoint = $d.Range($pStart, $pStart)
$insPoint.InsertAfter("https://automationstepbystep.com/")

# --- Step 2: delete everything from the end of that paragraph's text
# through the end of the document body (the "Thankyou..." paragraph, the
# four bullet/video paragraphs, and the trailing blank paragraphs), while
# leaving the sectPr (and therefore the section properties) untouched. ---
$p3after = $d.Paragraphs.Item(3).Range
$deleteStart = $p3after.End
$deleteEnd = $d.Content.End
$tail = $d.Range($deleteStart, $deleteEnd)
$tail.Delete()
